# 案件情報.xlsx — append newest scrape run (2025-10-24 01:14 JST)
#
# Two brand-new postings showed up at the top of the feed, and one more
# slipped in near the bottom (just before the last, lowest-priority row).
# Every row's "取得日時" (fetched-at) timestamp is refreshed to the new
# run time, since the whole sheet was re-scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = '2025-10-24 01:14:04'

# --- 1. Make room: insert a fresh row right under the header (row 2),
#        pushing the previous #1 listing (and everyone below it) down one.
$ws.Rows.Item(2).Insert()

# --- 2. Insert a second fresh row ahead of the final row. Before this
#        insert the old last row ("HPの微修正...") sits at row 15;
#        after it, that row becomes row 16 and row 15 is free for the
#        new listing that belongs there.
$ws.Rows.Item(15).Insert()

# --- 3. New listing #1 -> row 2 (top of the list)
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = '製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5419380'
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 298
$ws.Range("H2").Value = '🔥AI,Ai'

# --- 4. New listing #2 -> row 15 (inserted ahead of the final row)
$ws.Range("A15").Value = $newTimestamp
$ws.Range("B15").Value = '【急募】Google Play Consoleでのクローズテスト実施者募集!'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '~ 5,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5419425'
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("G15").Value = 10

# --- 5. Every pre-existing listing keeps its own data, but the scrape
#        timestamp in column A is refreshed for all of them too (rows
#        3-14 are the shifted originals #2-#13, row 16 is shifted
#        original #14).
$existingRows = @(3,4,5,6,7,8,9,10,11,12,13,14,16)
foreach ($r in $existingRows) {
    $ws.Range("A" + $r).Value = $newTimestamp
}

# --- 6. Hyperlinks: row inserts don't shift the sheet's stored hyperlink
#        anchors/relationships in this engine, so rebuild the collection
#        from scratch in final top-to-bottom order (F2..F16) to get
#        correct anchors and relationship targets.
$ws.Cells.Item(2, 6).Hyperlinks.Delete()

$linkUrls = @(
    'https://www.lancers.jp/work/detail/5419380',
    'https://www.lancers.jp/work/detail/5419191',
    'https://www.lancers.jp/work/detail/5418643',
    'https://www.lancers.jp/work/detail/5419221',
    'https://www.lancers.jp/work/detail/5418447',
    'https://www.lancers.jp/work/detail/5419226',
    'https://www.lancers.jp/work/detail/5418455',
    'https://www.lancers.jp/work/detail/5417544',
    'https://www.lancers.jp/work/detail/5418891',
    'https://www.lancers.jp/work/detail/5418644',
    'https://www.lancers.jp/work/detail/5418759',
    'https://www.lancers.jp/work/detail/5418738',
    'https://www.lancers.jp/work/detail/5418443',
    'https://www.lancers.jp/work/detail/5419425',
    'https://www.lancers.jp/work/detail/5418445'
)

for ($i = 0; $i -lt $linkUrls.Count; $i++) {
    $r = $i + 2
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $linkUrls[$i])
    $cell.Style = "Hyperlink"
}
